# Weekly price-sheet update: a new observation is inserted at the top of the
# "Cultivar IV Región" block (row 277), pushing all later rows down by one.
# The previously-last row (324) now lives at row 325, and the sheet's used
# range grows from A1:R324 to A1:R325.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh blank row at 277; Excel shifts rows 277-324 down to 278-325
# and carries the row-above's number formatting (style "2" on column D) with
# it, matching the target file's <c r="D277" s="2" .../> styling.
$ws.Rows(277).Insert()

$r = 277
$ws.Cells.Item($r, 1).Value  = 10
$ws.Cells.Item($r, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item($r, 3).Value  = "La Araucanía"
$ws.Cells.Item($r, 4).Value  = 45015
$ws.Cells.Item($r, 5).Value  = 9
$ws.Cells.Item($r, 6).Value  = 100112043
$ws.Cells.Item($r, 7).Value  = "Pepino dulce"
$ws.Cells.Item($r, 8).Value  = "Cultivar IV Región"
$ws.Cells.Item($r, 9).Value  = "Primera"
$ws.Cells.Item($r, 10).Value = 260
$ws.Cells.Item($r, 11).Value = 15000
$ws.Cells.Item($r, 12).Value = 16000
$ws.Cells.Item($r, 13).Value = 15538
$ws.Cells.Item($r, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item($r, 15).Value = "Provincia de Limarí"
$ws.Cells.Item($r, 16).Value = 863
$ws.Cells.Item($r, 17).Value = 18
$ws.Cells.Item($r, 18).Value = "Hortaliza"
